$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62 (G62 = 27781)
$ws.Cells.Item(62, 8).Value = 3353.348  # H62: 3403.2727 -> 3353.348
$ws.Cells.Item(62, 9).Value = 3095.4666  # I62: 3155.5 -> 3095.4666
$ws.Cells.Item(62, 11).Value = 3095.4666  # K62: 3155.5 -> 3095.4666
$ws.Cells.Item(62, 13).Value = -2471.4666  # M62: -2531.5 -> -2471.4666
# Row 65 (G65 = 27781)
$ws.Cells.Item(65, 8).Value = 3353.348  # H65: 3403.2727 -> 3353.348
$ws.Cells.Item(65, 9).Value = 3095.4666  # I65: 3155.5 -> 3095.4666
$ws.Cells.Item(65, 11).Value = 15477.333  # K65: 15777.5 -> 15477.333
$ws.Cells.Item(65, 13).Value = -12357.333  # M65: -12657.5 -> -12357.333
# Row 115 (G115 = 27957)
$ws.Cells.Item(115, 8).Value = 6250815  # H115: 6667552.5 -> 6250815
$ws.Cells.Item(115, 9).Value = 7692687.5  # I115: 8333716 -> 7692687.5
$ws.Cells.Item(115, 10).Value = 2700  # J115: 2900 -> 2700
$ws.Cells.Item(115, 11).Value = 23078062.5  # K115: 25001148 -> 23078062.5
$ws.Cells.Item(115, 12).Value = 8100  # L115: 8700 -> 8100
$ws.Cells.Item(115, 13).Value = -23076495.5  # M115: -24999581 -> -23076495.5
$ws.Cells.Item(115, 14).Value = -11234  # N115: -11834 -> -11234

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G2 = 27713)
$ws.Cells.Item(2, 8).Value = 2128.2778  # H2: 2176.8823 -> 2128.2778
$ws.Cells.Item(2, 9).Value = 2189.2  # I2: 2145.5454 -> 2189.2
$ws.Cells.Item(2, 10).Value = 2052.125  # J2: 2234.3333 -> 2052.125
$ws.Cells.Item(2, 11).Value = 2189.2  # K2: 2145.5454 -> 2189.2
$ws.Cells.Item(2, 12).Value = 2052.125  # L2: 2234.3333 -> 2052.125
$ws.Cells.Item(2, 13).Value = -2076.2  # M2: -2032.5454 -> -2076.2
$ws.Cells.Item(2, 14).Value = -2278.125  # N2: -2460.3333 -> -2278.125
# Row 32 (G32 = 44147)
$ws.Cells.Item(32, 8).Value = 14497399  # H32: 12824651 -> 14497399
$ws.Cells.Item(32, 9).Value = 15154963  # I32: 13336397 -> 15154963
$ws.Cells.Item(32, 11).Value = 15154963  # K32: 13336397 -> 15154963
$ws.Cells.Item(32, 13).Value = -15154676  # M32: -13336110 -> -15154676
# Row 39 (G39 = 2257)
$ws.Cells.Item(39, 8).Value = 4940  # H39: 3000 -> 4940
$ws.Cells.Item(39, 9).Value = 4940  # I39: 3000 -> 4940
$ws.Cells.Item(39, 11).Value = 4940  # K39: 3000 -> 4940
$ws.Cells.Item(39, 13).Value = -4420  # M39: -2480 -> -4420
# Row 61 (G61 = 43999)
$ws.Cells.Item(61, 8).Value = 879.5806  # H61: 913.6896400000001 -> 879.5806
$ws.Cells.Item(61, 9).Value = 807.2222  # I61: 841 -> 807.2222
$ws.Cells.Item(61, 11).Value = 807.2222  # K61: 841 -> 807.2222
$ws.Cells.Item(61, 13).Value = -595.2222  # M61: -629 -> -595.2222
# Row 74 (G74 = 44000)
$ws.Cells.Item(74, 8).Value = 1139.758  # H74: 1007.1781 -> 1139.758
$ws.Cells.Item(74, 9).Value = 1090.1296  # I74: 942.02985 -> 1090.1296
$ws.Cells.Item(74, 10).Value = 1474.75  # J74: 1734.6666 -> 1474.75
$ws.Cells.Item(74, 11).Value = 1090.1296  # K74: 942.02985 -> 1090.1296
$ws.Cells.Item(74, 12).Value = 1474.75  # L74: 1734.6666 -> 1474.75
$ws.Cells.Item(74, 13).Value = -216.1296  # M74: -68.02985000000001 -> -216.1296
$ws.Cells.Item(74, 14).Value = -3222.75  # N74: -3482.6666 -> -3222.75
# Row 77 (G77 = 44000)
$ws.Cells.Item(77, 8).Value = 1139.758  # H77: 1007.1781 -> 1139.758
$ws.Cells.Item(77, 9).Value = 1090.1296  # I77: 942.02985 -> 1090.1296
$ws.Cells.Item(77, 10).Value = 1474.75  # J77: 1734.6666 -> 1474.75
$ws.Cells.Item(77, 11).Value = 5450.648  # K77: 4710.14925 -> 5450.648
$ws.Cells.Item(77, 12).Value = 7373.75  # L77: 8673.333000000001 -> 7373.75
$ws.Cells.Item(77, 13).Value = -1082.648  # M77: -342.1492500000004 -> -1082.648
$ws.Cells.Item(77, 14).Value = -16109.75  # N77: -17409.333 -> -16109.75
# Row 96 (G96 = 18207)
$ws.Cells.Item(96, 8).Value = 23778  # H96: 0 -> 23778
$ws.Cells.Item(96, 10).Value = 23778  # J96: 0 -> 23778
$ws.Cells.Item(96, 12).Value = 23778  # L96: 0 -> 23778
$ws.Cells.Item(96, 14).Value = -29270  # N96: None -> -29270
# Row 102 (G102 = 19945)
$ws.Cells.Item(102, 8).Value = 2087.2856  # H102: 2183.5 -> 2087.2856
$ws.Cells.Item(102, 9).Value = 2000  # I102: 1997.5 -> 2000
$ws.Cells.Item(102, 10).Value = 2611  # J102: 2555.5 -> 2611
$ws.Cells.Item(102, 11).Value = 2000  # K102: 1997.5 -> 2000
$ws.Cells.Item(102, 12).Value = 2611  # L102: 2555.5 -> 2611
$ws.Cells.Item(102, 13).Value = -378  # M102: -375.5 -> -378
$ws.Cells.Item(102, 14).Value = -5855  # N102: -5799.5 -> -5855
# Row 116 (G116 = 27713)
$ws.Cells.Item(116, 8).Value = 2128.2778  # H116: 2176.8823 -> 2128.2778
$ws.Cells.Item(116, 9).Value = 2189.2  # I116: 2145.5454 -> 2189.2
$ws.Cells.Item(116, 10).Value = 2052.125  # J116: 2234.3333 -> 2052.125
$ws.Cells.Item(116, 11).Value = 2189.2  # K116: 2145.5454 -> 2189.2
$ws.Cells.Item(116, 12).Value = 2052.125  # L116: 2234.3333 -> 2052.125
$ws.Cells.Item(116, 13).Value = 104.8000000000002  # M116: 148.4546 -> 104.8000000000002
$ws.Cells.Item(116, 14).Value = -6640.125  # N116: -6822.3333 -> -6640.125
# Row 136 (G136 = 43999)
$ws.Cells.Item(136, 8).Value = 879.5806  # H136: 913.6896400000001 -> 879.5806
$ws.Cells.Item(136, 9).Value = 807.2222  # I136: 841 -> 807.2222
$ws.Cells.Item(136, 11).Value = 2421.6666  # K136: 2523 -> 2421.6666
$ws.Cells.Item(136, 13).Value = 128.3334  # M136: 27 -> 128.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G3 = 27713)
$ws.Cells.Item(3, 8).Value = 2128.2778  # H3: 2176.8823 -> 2128.2778
$ws.Cells.Item(3, 9).Value = 2189.2  # I3: 2145.5454 -> 2189.2
$ws.Cells.Item(3, 10).Value = 2052.125  # J3: 2234.3333 -> 2052.125
$ws.Cells.Item(3, 11).Value = 2189.2  # K3: 2145.5454 -> 2189.2
$ws.Cells.Item(3, 12).Value = 2052.125  # L3: 2234.3333 -> 2052.125
$ws.Cells.Item(3, 13).Value = -2075.2  # M3: -2031.5454 -> -2075.2
$ws.Cells.Item(3, 14).Value = -2280.125  # N3: -2462.3333 -> -2280.125
# Row 132 (G132 = 41855)
$ws.Cells.Item(132, 8).Value = 42000  # H132: 40884.617 -> 42000
$ws.Cells.Item(132, 10).Value = 42000  # J132: 40884.617 -> 42000
$ws.Cells.Item(132, 12).Value = 42000  # L132: 40884.617 -> 42000
$ws.Cells.Item(132, 14).Value = -52120  # N132: -51004.617 -> -52120
# Row 134 (G134 = 43998)
$ws.Cells.Item(134, 8).Value = 1541.1968  # H134: 1486 -> 1541.1968
$ws.Cells.Item(134, 9).Value = 1187.6305  # I134: 1143.2858 -> 1187.6305
$ws.Cells.Item(134, 10).Value = 2625.4666  # J134: 2685.5 -> 2625.4666
$ws.Cells.Item(134, 11).Value = 3562.8915  # K134: 3429.8574 -> 3562.8915
$ws.Cells.Item(134, 12).Value = 7876.399800000001  # L134: 8056.5 -> 7876.399800000001
$ws.Cells.Item(134, 13).Value = -1027.8915  # M134: -894.8574000000003 -> -1027.8915
$ws.Cells.Item(134, 14).Value = -12946.3998  # N134: -13126.5 -> -12946.3998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 26 (G26 = 2004)
$ws.Cells.Item(26, 8).Value = 0  # H26: 10000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # J26: 10000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # L26: 10000 -> 0
$ws.Cells.Item(26, 14).Value = ""  # N26: clear (was -10574)
# Row 31 (G31 = 44023)
$ws.Cells.Item(31, 8).Value = 2138.6826  # H31: 2383.4546 -> 2138.6826
$ws.Cells.Item(31, 9).Value = 1241.7174  # I31: 1409.5385 -> 1241.7174
$ws.Cells.Item(31, 10).Value = 4565.7646  # J31: 4757.375 -> 4565.7646
$ws.Cells.Item(31, 11).Value = 1241.7174  # K31: 1409.5385 -> 1241.7174
$ws.Cells.Item(31, 12).Value = 4565.7646  # L31: 4757.375 -> 4565.7646
$ws.Cells.Item(31, 13).Value = -946.7174  # M31: -1114.5385 -> -946.7174
$ws.Cells.Item(31, 14).Value = -5155.7646  # N31: -5347.375 -> -5155.7646
# Row 34 (G34 = 44023)
$ws.Cells.Item(34, 8).Value = 2138.6826  # H34: 2383.4546 -> 2138.6826
$ws.Cells.Item(34, 9).Value = 1241.7174  # I34: 1409.5385 -> 1241.7174
$ws.Cells.Item(34, 10).Value = 4565.7646  # J34: 4757.375 -> 4565.7646
$ws.Cells.Item(34, 11).Value = 1241.7174  # K34: 1409.5385 -> 1241.7174
$ws.Cells.Item(34, 12).Value = 4565.7646  # L34: 4757.375 -> 4565.7646
$ws.Cells.Item(34, 13).Value = -1039.7174  # M34: -1207.5385 -> -1039.7174
$ws.Cells.Item(34, 14).Value = -4969.7646  # N34: -5161.375 -> -4969.7646

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113 (G113 = 27710)
$ws.Cells.Item(113, 8).Value = 11149.454  # H113: 9285 -> 11149.454
$ws.Cells.Item(113, 9).Value = 1506.2  # I113: 798.8570999999999 -> 1506.2
$ws.Cells.Item(113, 11).Value = 1506.2  # K113: 798.8570999999999 -> 1506.2
$ws.Cells.Item(113, 13).Value = 663.8  # M113: 1371.1429 -> 663.8
# Row 122 (G122 = 36182)
$ws.Cells.Item(122, 8).Value = 9092634  # H122: 25013650 -> 9092634
$ws.Cells.Item(122, 9).Value = 11112553  # I122: 25013650 -> 11112553
$ws.Cells.Item(122, 10).Value = 3000  # J122: 0 -> 3000
$ws.Cells.Item(122, 11).Value = 33337659  # K122: 75040950 -> 33337659
$ws.Cells.Item(122, 12).Value = 9000  # L122: 0 -> 9000
$ws.Cells.Item(122, 13).Value = -33335209  # M122: -75038500 -> -33335209
$ws.Cells.Item(122, 14).Value = -13900  # N122: None -> -13900
# Row 132 (G132 = 44008)
$ws.Cells.Item(132, 8).Value = 3577.6743  # H132: 3455.5107 -> 3577.6743
$ws.Cells.Item(132, 9).Value = 3585.6924  # I132: 3669.2632 -> 3585.6924
$ws.Cells.Item(132, 10).Value = 3499.5  # J132: 2553 -> 3499.5
$ws.Cells.Item(132, 11).Value = 10757.0772  # K132: 11007.7896 -> 10757.0772
$ws.Cells.Item(132, 12).Value = 10498.5  # L132: 7659 -> 10498.5
$ws.Cells.Item(132, 13).Value = -8227.0772  # M132: -8477.7896 -> -8227.0772
$ws.Cells.Item(132, 14).Value = -15558.5  # N132: -12719 -> -15558.5
# Row 138 (G138 = 42325)
$ws.Cells.Item(138, 8).Value = 23329  # H138: 25395.666 -> 23329
$ws.Cells.Item(138, 10).Value = 23329  # J138: 25395.666 -> 23329
$ws.Cells.Item(138, 12).Value = 23329  # L138: 25395.666 -> 23329
$ws.Cells.Item(138, 14).Value = -33609  # N138: -35675.666 -> -33609

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G7 = 36249)
$ws.Cells.Item(7, 8).Value = 2341.524  # H7: 2357.3333 -> 2341.524
$ws.Cells.Item(7, 9).Value = 2357.4736  # I7: 2368.5217 -> 2357.4736
$ws.Cells.Item(7, 10).Value = 2190  # J7: 2100 -> 2190
$ws.Cells.Item(7, 11).Value = 2357.4736  # K7: 2368.5217 -> 2357.4736
$ws.Cells.Item(7, 12).Value = 2190  # L7: 2100 -> 2190
$ws.Cells.Item(7, 13).Value = -2245.4736  # M7: -2256.5217 -> -2245.4736
$ws.Cells.Item(7, 14).Value = -2414  # N7: -2324 -> -2414
# Row 126 (G126 = 36249)
$ws.Cells.Item(126, 8).Value = 2341.524  # H126: 2357.3333 -> 2341.524
$ws.Cells.Item(126, 9).Value = 2357.4736  # I126: 2368.5217 -> 2357.4736
$ws.Cells.Item(126, 10).Value = 2190  # J126: 2100 -> 2190
$ws.Cells.Item(126, 11).Value = 7072.4208  # K126: 7105.5651 -> 7072.4208
$ws.Cells.Item(126, 12).Value = 6570  # L126: 6300 -> 6570
$ws.Cells.Item(126, 13).Value = -4602.4208  # M126: -4635.5651 -> -4602.4208
$ws.Cells.Item(126, 14).Value = -11510  # N126: -11240 -> -11510
# Row 132 (G132 = 44058)
$ws.Cells.Item(132, 8).Value = 945.4  # H132: 1103.42 -> 945.4
$ws.Cells.Item(132, 9).Value = 944.2659  # I132: 1051.1686 -> 944.2659
$ws.Cells.Item(132, 10).Value = 963.1667  # J132: 1526.1818 -> 963.1667
$ws.Cells.Item(132, 11).Value = 2832.7977  # K132: 3153.5058 -> 2832.7977
$ws.Cells.Item(132, 12).Value = 2889.5001  # L132: 4578.5454 -> 2889.5001
$ws.Cells.Item(132, 13).Value = -302.7977000000001  # M132: -623.5057999999999 -> -302.7977000000001
$ws.Cells.Item(132, 14).Value = -7949.5001  # N132: -9638.545399999999 -> -7949.5001
# Row 137 (G137 = 43296)
$ws.Cells.Item(137, 8).Value = 35951.6  # H137: 28720.727 -> 35951.6
$ws.Cells.Item(137, 10).Value = 35951.6  # J137: 28720.727 -> 35951.6
$ws.Cells.Item(137, 12).Value = 35951.6  # L137: 28720.727 -> 35951.6
$ws.Cells.Item(137, 14).Value = -46151.6  # N137: -38920.727 -> -46151.6
# Row 139 (G139 = 43310)
$ws.Cells.Item(139, 8).Value = 0  # H139: 41750 -> 0
$ws.Cells.Item(139, 10).Value = 0  # J139: 41750 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 41750 -> 0
$ws.Cells.Item(139, 14).Value = ""  # N139: clear (was -52030)
# Row 141 (G141 = 42487)
$ws.Cells.Item(141, 8).Value = 60000  # H141: 65000 -> 60000
$ws.Cells.Item(141, 10).Value = 60000  # J141: 65000 -> 60000
$ws.Cells.Item(141, 12).Value = 60000  # L141: 65000 -> 60000
$ws.Cells.Item(141, 14).Value = -70360  # N141: -75360 -> -70360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41 (G41 = 21725)
$ws.Cells.Item(41, 8).Value = 3377  # H41: 6338.5 -> 3377
$ws.Cells.Item(41, 10).Value = 3377  # J41: 6338.5 -> 3377
$ws.Cells.Item(41, 12).Value = 3377  # L41: 6338.5 -> 3377
$ws.Cells.Item(41, 14).Value = -4157  # N41: -7118.5 -> -4157
# Row 45 (G45 = 21726)
$ws.Cells.Item(45, 8).Value = 11281.375  # H45: 12210.556 -> 11281.375
$ws.Cells.Item(45, 10).Value = 10321.571  # J45: 11486.875 -> 10321.571
$ws.Cells.Item(45, 12).Value = 10321.571  # L45: 11486.875 -> 10321.571
$ws.Cells.Item(45, 14).Value = -11303.571  # N45: -12468.875 -> -11303.571
# Row 49 (G49 = 3397)
$ws.Cells.Item(49, 8).Value = 20000  # H49: 5000 -> 20000
$ws.Cells.Item(49, 10).Value = 20000  # J49: 5000 -> 20000
$ws.Cells.Item(49, 12).Value = 20000  # L49: 5000 -> 20000
$ws.Cells.Item(49, 14).Value = -20460  # N49: -5460 -> -20460
# Row 107 (G107 = 27746)
$ws.Cells.Item(107, 8).Value = 304.66666  # H107: 282.66666 -> 304.66666
$ws.Cells.Item(107, 9).Value = 283  # I107: 249.83333 -> 283
$ws.Cells.Item(107, 11).Value = 849  # K107: 749.49999 -> 849
$ws.Cells.Item(107, 13).Value = 1071  # M107: 1170.50001 -> 1071
# Row 108 (G108 = 25661)
$ws.Cells.Item(108, 8).Value = 40000  # H108: 0 -> 40000
$ws.Cells.Item(108, 10).Value = 40000  # J108: 0 -> 40000
$ws.Cells.Item(108, 12).Value = 40000  # L108: 0 -> 40000
$ws.Cells.Item(108, 14).Value = -47680  # N108: None -> -47680
# Row 122 (G122 = 36208)
$ws.Cells.Item(122, 8).Value = 2093.125  # H122: 2478.8 -> 2093.125
$ws.Cells.Item(122, 9).Value = 1945.1  # I122: 2402.9333 -> 1945.1
$ws.Cells.Item(122, 10).Value = 2241.15  # J122: 2554.6667 -> 2241.15
$ws.Cells.Item(122, 11).Value = 5835.299999999999  # K122: 7208.7999 -> 5835.299999999999
$ws.Cells.Item(122, 12).Value = 6723.450000000001  # L122: 7664.000100000001 -> 6723.450000000001
$ws.Cells.Item(122, 13).Value = -3385.299999999999  # M122: -4758.7999 -> -3385.299999999999
$ws.Cells.Item(122, 14).Value = -11623.45  # N122: -12564.0001 -> -11623.45
